$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44628
$ws.Range("K2").Value = "Black Amber"
$ws.Range("M2").Value = 270
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("S2").Value = 861

# Row 3
$ws.Range("D3").Value = 44217
$ws.Range("K3").Value = "Black Amber"
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 16500
$ws.Range("Q3").Value = "$/bandeja 18 kilos granel"
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 917

# Row 4
$ws.Range("D4").Value = 44944
$ws.Range("K4").Value = "Larry Ann"
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 7500
$ws.Range("S4").Value = 417

# Row 5
$ws.Range("D5").Value = 44921
$ws.Range("K5").Value = "Angeleno"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 450
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19111
$ws.Range("S5").Value = 1062

# Row 7
$ws.Range("D7").Value = 44596
$ws.Range("M7").Value = 250

# Row 9
$ws.Range("D9").Value = 44229
$ws.Range("K9").Value = "Fortuna"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("Q9").Value = "$/bandeja 18 kilos granel"
$ws.Range("S9").Value = 806

# Row 10
$ws.Range("D10").Value = 44987
$ws.Range("M10").Value = 400
$ws.Range("N10").Value = 5000
$ws.Range("O10").Value = 6000
$ws.Range("P10").Value = 5750
$ws.Range("Q10").Value = "$/bandeja 10 kilos"
$ws.Range("R10").Value = "Región Metropolitana"
$ws.Range("S10").Value = 575
$ws.Range("T10").Value = 10

# Row 11
$ws.Range("D11").Value = 44953
$ws.Range("M11").Value = 350
$ws.Range("N11").Value = 19000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 19571
$ws.Range("S11").Value = 1087

# Row 12
$ws.Range("D12").Value = 44614
$ws.Range("M12").Value = 250
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 19000
$ws.Range("P12").Value = 18500
$ws.Range("Q12").Value = "$/bandeja 18 kilos granel"
$ws.Range("S12").Value = 1028
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44999
$ws.Range("N13").Value = 19000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 19500
$ws.Range("S13").Value = 1083

# Row 14
$ws.Range("D14").Value = 44973
$ws.Range("K14").Value = "Larry Ann"
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 18000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 19000
$ws.Range("S14").Value = 1056

# Row 15
$ws.Range("D15").Value = 44973
$ws.Range("K15").Value = "Pink Delight"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 270
$ws.Range("N15").Value = 19000
$ws.Range("P15").Value = 19500
$ws.Range("S15").Value = 1083

# Row 16
$ws.Range("D16").Value = 44243
$ws.Range("K16").Value = "Black Amber"
$ws.Range("L16").Value = "Primera"
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 14500
$ws.Range("Q16").Value = "$/caja 18 kilos granel"
$ws.Range("S16").Value = 806

# Row 17
$ws.Range("D17").Value = 44175
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 21000
$ws.Range("O17").Value = 22000
$ws.Range("P17").Value = 21500
$ws.Range("S17").Value = 1194

# Row 18
$ws.Range("D18").Value = 45021
$ws.Range("K18").Value = "Angeleno"
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 300

# Row 19
$ws.Range("D19").Value = 45013
$ws.Range("M19").Value = 280
$ws.Range("N19").Value = 21000
$ws.Range("O19").Value = 22000
$ws.Range("P19").Value = 21536
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 1196

# Row 20
$ws.Range("D20").Value = 45008
$ws.Range("K20").Value = "Angeleno"
$ws.Range("N20").Value = 19000
$ws.Range("P20").Value = 19500
$ws.Range("S20").Value = 1083

# Row 21
$ws.Range("D21").Value = 45008
$ws.Range("K21").Value = "Fortuna"
$ws.Range("L21").Value = "Tercera"

# Row 22
$ws.Range("D22").Value = 44952
$ws.Range("K22").Value = "Larry Ann"
$ws.Range("L22").Value = "Primera"
$ws.Range("N22").Value = 20000
$ws.Range("O22").Value = 21000
$ws.Range("P22").Value = 20500
$ws.Range("R22").Value = "Provincia de Curicó"
$ws.Range("S22").Value = 1139

# Row 23
$ws.Range("D23").Value = 44587
$ws.Range("M23").Value = 300
$ws.Range("N23").Value = 15000
$ws.Range("O23").Value = 16000
$ws.Range("P23").Value = 15500
$ws.Range("Q23").Value = "$/caja 18 kilos granel"
$ws.Range("R23").Value = "Región de O'Higgins"
$ws.Range("S23").Value = 861

# Row 26
$ws.Range("D26").Value = 44580
$ws.Range("K26").Value = "Black Amber"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 270
$ws.Range("N26").Value = 19000
$ws.Range("P26").Value = 19500
$ws.Range("Q26").Value = "$/bandeja 18 kilos granel"
$ws.Range("R26").Value = "Región Metropolitana"
$ws.Range("S26").Value = 1083

# Row 27
$ws.Range("D27").Value = 44994
$ws.Range("K27").Value = "Angeleno"
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 270
$ws.Range("N27").Value = 19000
$ws.Range("O27").Value = 20000
$ws.Range("P27").Value = 19500
$ws.Range("R27").Value = "Región Metropolitana"
$ws.Range("S27").Value = 1083

# Row 28
$ws.Range("D28").Value = 44966
$ws.Range("K28").Value = "Sapphire"
$ws.Range("M28").Value = 250
$ws.Range("N28").Value = 18000
$ws.Range("O28").Value = 20000
$ws.Range("P28").Value = 19000
$ws.Range("S28").Value = 1056

# Row 29
$ws.Range("D29").Value = 44574
$ws.Range("K29").Value = "Black Amber"
$ws.Range("L29").Value = "Primera"
$ws.Range("N29").Value = 18000
$ws.Range("O29").Value = 19000
$ws.Range("P29").Value = 18500
$ws.Range("S29").Value = 1028

# Row 30
$ws.Range("D30").Value = 44245
$ws.Range("K30").Value = "Black Amber"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 250
$ws.Range("N30").Value = 14000
$ws.Range("O30").Value = 15000
$ws.Range("P30").Value = 14500
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 806

# Row 31
$ws.Range("D31").Value = 44239
$ws.Range("K31").Value = "Fortuna"
$ws.Range("N31").Value = 15000
$ws.Range("O31").Value = 16000
$ws.Range("P31").Value = 15500
$ws.Range("S31").Value = 861

# Row 32
$ws.Range("D32").Value = 44174
$ws.Range("M32").Value = 270
$ws.Range("N32").Value = 20000
$ws.Range("O32").Value = 21000
$ws.Range("P32").Value = 20500
$ws.Range("Q32").Value = "$/caja 18 kilos granel"
$ws.Range("S32").Value = 1139

# Row 33
$ws.Range("D33").Value = 45030
$ws.Range("K33").Value = "Fortuna"
$ws.Range("L33").Value = "Segunda"
$ws.Range("M33").Value = 200
$ws.Range("N33").Value = 19000
$ws.Range("O33").Value = 20000
$ws.Range("P33").Value = 19500
$ws.Range("S33").Value = 1083

# Row 34
$ws.Range("D34").Value = 44285
$ws.Range("K34").Value = "Angeleno"
$ws.Range("N34").Value = 14000
$ws.Range("O34").Value = 15000
$ws.Range("P34").Value = 14500
$ws.Range("S34").Value = 806

# Row 35
$ws.Range("D35").Value = 44706
$ws.Range("K35").Value = "Angeleno"
$ws.Range("N35").Value = 15000
$ws.Range("O35").Value = 16000
$ws.Range("P35").Value = 15500
$ws.Range("S35").Value = 861

# Row 36
$ws.Range("D36").Value = 44650
$ws.Range("K36").Value = "Angeleno"
$ws.Range("N36").Value = 17000
$ws.Range("O36").Value = 18000
$ws.Range("P36").Value = 17500
$ws.Range("R36").Value = "Región de O'Higgins"
$ws.Range("S36").Value = 972

# Row 37
$ws.Range("D37").Value = 44169
$ws.Range("K37").Value = "Angeleno"
$ws.Range("L37").Value = "Tercera"
$ws.Range("M37").Value = 250
$ws.Range("N37").Value = 24000
$ws.Range("O37").Value = 25000
$ws.Range("P37").Value = 24500
$ws.Range("R37").Value = "Región de O'Higgins"
$ws.Range("S37").Value = 1361

# Row 38
$ws.Range("D38").Value = 45002
$ws.Range("L38").Value = "Segunda"
$ws.Range("M38").Value = 300
$ws.Range("N38").Value = 21000
$ws.Range("O38").Value = 22000
$ws.Range("P38").Value = 21500
$ws.Range("S38").Value = 1194
